# Auto-generated Excel COM-interop script
# Applies cached-value updates (market price snapshot refresh) to the
# Famfrit_Profits workbook sheets, as captured by the scheduled runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 244
$ws.Range("I11").Value = 244
$ws.Range("K11").Value = 244
$ws.Range("M11").Value = -104

$ws.Range("H33").Value = 17741.166
$ws.Range("I33").Value = 23527.223
$ws.Range("J33").Value = 383
$ws.Range("K33").Value = 23527.223
$ws.Range("L33").Value = 383
$ws.Range("M33").Value = -23298.223
$ws.Range("N33").Value = -841

$ws.Range("H74").Value = 5579.8125
$ws.Range("I74").Value = 5039.5713
$ws.Range("K74").Value = 5039.5713
$ws.Range("M74").Value = -4103.5713

$ws.Range("H77").Value = 5579.8125
$ws.Range("I77").Value = 5039.5713
$ws.Range("K77").Value = 25197.8565
$ws.Range("M77").Value = -20517.8565

$ws.Range("H92").Value = 1982.8
$ws.Range("I92").Value = 2293.3333
$ws.Range("J92").Value = 1517
$ws.Range("K92").Value = 2293.3333
$ws.Range("L92").Value = 1517
$ws.Range("M92").Value = -1045.3333
$ws.Range("N92").Value = -4013

$ws.Range("H132").Value = 14286825
$ws.Range("I132").Value = 14926409
$ws.Range("J132").Value = 2770.6667
$ws.Range("K132").Value = 44779227
$ws.Range("L132").Value = 8312.000100000001
$ws.Range("M132").Value = -44776697
$ws.Range("N132").Value = -13372.0001

$ws.Range("H137").Value = 16970.908
$ws.Range("I137").Value = 5248.125
$ws.Range("J137").Value = 48231.668
$ws.Range("K137").Value = 15744.375
$ws.Range("L137").Value = 144695.004
$ws.Range("M137").Value = -13194.375
$ws.Range("N137").Value = -149795.004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 257.5
$ws.Range("I17").Value = 121
$ws.Range("J17").Value = 394
$ws.Range("K17").Value = 121
$ws.Range("L17").Value = 394
$ws.Range("M17").Value = 52
$ws.Range("N17").Value = -740

$ws.Range("H32").Value = 4564.0366
$ws.Range("J32").Value = 7093.1113
$ws.Range("L32").Value = 7093.1113
$ws.Range("N32").Value = -7667.1113

$ws.Range("H45").Value = 2487.8572
$ws.Range("I45").Value = 1700
$ws.Range("J45").Value = 2803
$ws.Range("K45").Value = 1700
$ws.Range("L45").Value = 2803
$ws.Range("M45").Value = -1323
$ws.Range("N45").Value = -3557

$ws.Range("H61").Value = 37041570
$ws.Range("I61").Value = 45458310
$ws.Range("J61").Value = 7899.8
$ws.Range("K61").Value = 45458310
$ws.Range("L61").Value = 7899.8
$ws.Range("M61").Value = -45458098
$ws.Range("N61").Value = -8323.799999999999

$ws.Range("H74").Value = 27058320
$ws.Range("I74").Value = 29445648
$ws.Range("J74").Value = 1933
$ws.Range("K74").Value = 29445648
$ws.Range("L74").Value = 1933
$ws.Range("M74").Value = -29444774
$ws.Range("N74").Value = -3681

$ws.Range("H77").Value = 27058320
$ws.Range("I77").Value = 29445648
$ws.Range("J77").Value = 1933
$ws.Range("K77").Value = 147228240
$ws.Range("L77").Value = 9665
$ws.Range("M77").Value = -147223872
$ws.Range("N77").Value = -18401

$ws.Range("H132").Value = 18231232
$ws.Range("I132").Value = 6359.3403
$ws.Range("K132").Value = 19078.0209
$ws.Range("M132").Value = -16548.0209

$ws.Range("H133").Value = 119960
$ws.Range("J133").Value = 119960
$ws.Range("L133").Value = 119960
$ws.Range("N133").Value = -125020

$ws.Range("H136").Value = 37041570
$ws.Range("I136").Value = 45458310
$ws.Range("J136").Value = 7899.8
$ws.Range("K136").Value = 136374930
$ws.Range("L136").Value = 23699.4
$ws.Range("M136").Value = -136372380
$ws.Range("N136").Value = -28799.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 384874.5
$ws.Range("J95").Value = 384874.5
$ws.Range("L95").Value = 384874.5
$ws.Range("N95").Value = -390366.5

$ws.Range("H105").Value = 9505.857
$ws.Range("I105").Value = 13321.625
$ws.Range("J105").Value = 4418.1665
$ws.Range("K105").Value = 13321.625
$ws.Range("L105").Value = 4418.1665
$ws.Range("M105").Value = -11574.625
$ws.Range("N105").Value = -7912.1665

$ws.Range("H134").Value = 2858887.8
$ws.Range("I134").Value = 3031911.2
$ws.Range("K134").Value = 9095733.600000001
$ws.Range("M134").Value = -9093198.600000001

$ws.Range("H139").Value = 199999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 199999
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 199999
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -210279

$ws.Range("H141").Value = 69333.336
$ws.Range("I141").Value = 65000
$ws.Range("J141").Value = 78000
$ws.Range("K141").Value = 65000
$ws.Range("L141").Value = 78000
$ws.Range("M141").Value = -59820
$ws.Range("N141").Value = -88360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 355.36365
$ws.Range("I10").Value = 367
$ws.Range("J10").Value = 335
$ws.Range("K10").Value = 367
$ws.Range("L10").Value = 335
$ws.Range("M10").Value = -228
$ws.Range("N10").Value = -613

$ws.Range("H16").Value = 1507.5
$ws.Range("I16").Value = 1217.8
$ws.Range("K16").Value = 1217.8
$ws.Range("M16").Value = -930.8

$ws.Range("H22").Value = 7022.533
$ws.Range("J22").Value = 424.5
$ws.Range("L22").Value = 424.5
$ws.Range("N22").Value = -1124.5

$ws.Range("H31").Value = 29415814
$ws.Range("I31").Value = 2490.1875
$ws.Range("J31").Value = 55560988
$ws.Range("K31").Value = 2490.1875
$ws.Range("L31").Value = 55560988
$ws.Range("M31").Value = -2195.1875
$ws.Range("N31").Value = -55561578

$ws.Range("H34").Value = 29415814
$ws.Range("I34").Value = 2490.1875
$ws.Range("J34").Value = 55560988
$ws.Range("K34").Value = 2490.1875
$ws.Range("L34").Value = 55560988
$ws.Range("M34").Value = -2288.1875
$ws.Range("N34").Value = -55561392

$ws.Range("H99").Value = 6119.2
$ws.Range("I99").Value = 6114.8887
$ws.Range("J99").Value = 6130.2856
$ws.Range("K99").Value = 6114.8887
$ws.Range("L99").Value = 6130.2856
$ws.Range("M99").Value = -4616.8887
$ws.Range("N99").Value = -9126.285599999999

$ws.Range("H113").Value = 1507.5
$ws.Range("I113").Value = 1217.8
$ws.Range("K113").Value = 1217.8
$ws.Range("M113").Value = 952.2

$ws.Range("H126").Value = 6119.2
$ws.Range("I126").Value = 6114.8887
$ws.Range("J126").Value = 6130.2856
$ws.Range("K126").Value = 18344.6661
$ws.Range("L126").Value = 18390.8568
$ws.Range("M126").Value = -15874.6661
$ws.Range("N126").Value = -23330.8568

$ws.Range("H131").Value = 22006.076
$ws.Range("I131").Value = 12999
$ws.Range("J131").Value = 23643.727
$ws.Range("K131").Value = 12999
$ws.Range("L131").Value = 23643.727
$ws.Range("M131").Value = -7959
$ws.Range("N131").Value = -33723.727

$ws.Range("H132").Value = 87153.52
$ws.Range("I132").Value = 100350.586
$ws.Range("K132").Value = 301051.758
$ws.Range("M132").Value = -298521.758

$ws.Range("H141").Value = 107327.695
$ws.Range("J141").Value = 107327.695
$ws.Range("L141").Value = 107327.695
$ws.Range("N141").Value = -117687.695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1085.6
$ws.Range("I107").Value = 611.2857
$ws.Range("K107").Value = 1833.8571
$ws.Range("M107").Value = 86.14289999999983

$ws.Range("H113").Value = 1842.7916
$ws.Range("J113").Value = 2178.2222
$ws.Range("L113").Value = 6534.6666
$ws.Range("N113").Value = -10874.6666

$ws.Range("H128").Value = 99450
$ws.Range("I128").Value = 99450
$ws.Range("K128").Value = 298350
$ws.Range("M128").Value = -293370

$ws.Range("H139").Value = 2385.2666
$ws.Range("I139").Value = 2285.6155
$ws.Range("J139").Value = 3033
$ws.Range("K139").Value = 6856.8465
$ws.Range("L139").Value = 9099
$ws.Range("M139").Value = -1716.8465
$ws.Range("N139").Value = -19379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3276.85
$ws.Range("I102").Value = 2588.1667
$ws.Range("K102").Value = 2588.1667
$ws.Range("M102").Value = -966.1667000000002

$ws.Range("H113").Value = 3847.1714
$ws.Range("I113").Value = 3090.4707
$ws.Range("J113").Value = 4561.8335
$ws.Range("K113").Value = 3090.4707
$ws.Range("L113").Value = 4561.8335
$ws.Range("M113").Value = -920.4706999999999
$ws.Range("N113").Value = -8901.833500000001

$ws.Range("H132").Value = 6317.8423
$ws.Range("I132").Value = 4651.625
$ws.Range("K132").Value = 13954.875
$ws.Range("M132").Value = -11424.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2624.5557
$ws.Range("I40").Value = 2590.12
$ws.Range("J40").Value = 3055
$ws.Range("K40").Value = 2590.12
$ws.Range("L40").Value = 3055
$ws.Range("M40").Value = -2454.12
$ws.Range("N40").Value = -3327

$ws.Range("H46").Value = 2268.1333
$ws.Range("I46").Value = 831.5789
$ws.Range("K46").Value = 831.5789
$ws.Range("M46").Value = -643.5789

$ws.Range("H122").Value = 3155.1277
$ws.Range("I122").Value = 2596.3333
$ws.Range("J122").Value = 3738.2173
$ws.Range("K122").Value = 7788.999899999999
$ws.Range("L122").Value = 11214.6519
$ws.Range("M122").Value = -5338.999899999999
$ws.Range("N122").Value = -16114.6519

$ws.Range("H131").Value = 77888
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 77888
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 77888
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -87968

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 45909440
$ws.Range("I100").Value = 50500268
$ws.Range("K100").Value = 101000536
$ws.Range("M100").Value = -100999995

$ws.Range("H136").Value = 2178.4722
$ws.Range("I136").Value = 773.3103599999999
$ws.Range("J136").Value = 7999.857
$ws.Range("K136").Value = 2319.93108
$ws.Range("L136").Value = 23999.571
$ws.Range("M136").Value = 230.0689200000002
$ws.Range("N136").Value = -29099.571

